$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 holds the column headers for the q_acct table. Two new columns are
# being added: "Account" right after id/q_acct, and "active" right before
# owner_id. Shift the existing headers right (starting from the rightmost
# cell) to make room, then fill in the two new header cells.

$ws.Range("K5").Value = $ws.Range("I5").Text    # owner_id   -> K5
$ws.Range("J5").Value = "active"                 # new header -> J5
$ws.Range("I5").Value = $ws.Range("H5").Text    # attachment -> I5
$ws.Range("H5").Value = $ws.Range("G5").Text    # Balance    -> H5
$ws.Range("G5").Value = $ws.Range("F5").Text    # Debit      -> G5
$ws.Range("F5").Value = $ws.Range("E5").Text    # Credit     -> F5
$ws.Range("E5").Value = $ws.Range("D5").Text    # Description-> E5
$ws.Range("D5").Value = $ws.Range("C5").Text    # Date Entered -> D5
$ws.Range("C5").Value = "Account"                # new header -> C5

$null = $ws.Range("K6").Select()
